$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.272.14'
$ws.Range("E2").Value = '  +1.79%  '
$ws.Range("D3").Value = '3.445.44'
$ws.Range("E3").Value = '  +2.19%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''414.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.08%  '
$ws.Range("D6").Value = '''130.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.65%  '
$ws.Range("D7").Value = '''0.630'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6.51%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").Value = '''0.755'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +12.34%  '
$ws.Range("D10").Value = '''0.140'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +16.94%  '
$ws.Range("D11").Value = '''43.43'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.26%  '
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").Value = '''8.99'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.61%  '
$ws.Range("D14").Value = '''20.72'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.85%  '
$ws.Range("D15").Value = '''0.0000198'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +56.34%  '
$ws.Range("B16").Value = 'Uniswap'
$ws.Range("C16").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D16").Value = '''13.23'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +19.95%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.457.21'
$ws.Range("E17").Value = '  +3.49%  '
$ws.Range("E18").Value = '  +4.80%  '
$ws.Range("D19").Value = '62.241.70'
$ws.Range("D20").Value = '''399.68'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +27.05%  '
$ws.Range("D21").Value = '''90.57'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +8.34%  '
$ws.Range("D22").Value = '''3.21'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("E23").Value = '  +4.86%  '
$ws.Range("D24").Value = '''3.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.22%  '
$ws.Range("D25").Value = '''34.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +16.30%  '
$ws.Range("B26").Value = 'Filecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D26").Value = '''8.60'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.15%  '
$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").Value = '''4.80'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.58%  '
$ws.Range("D28").Value = '''7.67'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.40%  '
$ws.Range("D29").Value = '''2.74'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +9.93%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = '''0.117'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.60%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = '''44.15'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.88%  '
$ws.Range("E32").Value = '  -0.40%  '
$ws.Range("D33").Value = '''11.91'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.38%  '
$ws.Range("D34").Value = '''0.999'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("D35").Value = '''0.0497'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.59%  '
$ws.Range("D36").Value = '''52.48'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.65%  '
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("E38").Value = '  -0.29%  '
$ws.Range("D39").Value = '''2.91'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.56%  '
$ws.Range("E40").Value = '  +6.95%  '
$ws.Range("E41").Value = '  +8.03%  '
$ws.Range("D42").Value = '''140.52'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.03%  '
$ws.Range("D43").Value = '''2.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.06%  '
$ws.Range("D44").Value = '''4.06'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.23%  '
$ws.Range("D45").Value = '''16.90'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.49%  '
$ws.Range("D46").Value = '''2.34'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.18%  '
$ws.Range("D47").Value = '''22.62'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.19%  '
$ws.Range("D48").Value = '2.129.59'
$ws.Range("E48").Value = '  +0.31%  '
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").Value = '''2.28'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.19%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").Value = '''1.94'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.01%  '
$ws.Range("D51").Value = '''0.0371'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +9.08%  '
